$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the description of the existing "G forces" row (row 10) to
# "Gs of Acceleration" since it now specifically refers to acceleration.
$ws.Range("A10").Value = "Gs of Acceleration"

# Add a new row (row 11) describing "Gs of Deceleration" with "Dx" as the
# variable name, matching columns B (var_name) and C (unit).
$ws.Range("A11").Value = "Gs of Deceleration"
$ws.Range("B11").Value = "Dx"
$ws.Range("C11").Value = "Gs"

# Move the active selection below the newly added data, mirroring the
# author's final cursor placement.
$ws.Range("A13").Select()
